# The document has a single section whose primary ("default") and
# first-page headers/footers each contain one inline picture:
#   - Headers: the BTec_Logo-Orange picture, currently named "image1.jpg",
#     should be renamed to "image2.jpg".
#   - Footers: the Pearson logo picture, currently named "image2.png",
#     should be renamed to "image1.png".
#
# Word's InlineShapes don't carry visible text, so Find/Replace can't
# reach them - update the picture's Name property directly via the
# InlineShape object exposed on each header/footer Range.

$d = $word.ActiveDocument
$section = $d.Sections(1)

# --- Headers: BTec logo image1.jpg -> image2.jpg ---
for ($i = 1; $i -le 3; $i++) {
    $header = $section.Headers($i)
    if ($header.Exists) {
        $shapes = $header.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shapes.Item($j).Name = "image2.jpg"
        }
    }
}

# --- Footers: Pearson logo image2.png -> image1.png ---
for ($i = 1; $i -le 3; $i++) {
    $footer = $section.Footers($i)
    if ($footer.Exists) {
        $shapes = $footer.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shapes.Item($j).Name = "image1.png"
        }
    }
}
